$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "#PBM:DWP:Defect:SmartIT/MidTier integration issue;  DWP not displaying customer requests.  Requires investigation.:R&D"
$ws.Range("F3").Value = "#PBM:DWP:Defect:Internal server error when opening delegate approvals for approvers with deleted profiles;  handling of deleted approver profiles in approval workflow needs to be improved.:R&D"
$ws.Range("F4").Value = "#PBM:DWPA:Entitlement:Entitlement not properly assigned to users; corrected entitlement configuration.:NA"
$ws.Range("F5").Value = "#PBM:DWPC:Defect:Dev DWPC forms: Missing data tables in action mapping fixed.:R&D"
$ws.Range("F6").Value = "#PBM:DWP:Workflow:Status update failure in Catalog; workflow investigated and corrected.:R&D"
$ws.Range("F7").Value = "#PBM:DWP:Upgrade:23.3.04 patch deployed to QA environment excluding sample data, DWPC catalogs, and DWP Studio pages.:NA"
$ws.Range("F8").Value = "#PBM:DWP:Integration:Investigating DWP-ITSM integration issue; worklog comments not syncing.:R&D"
$ws.Range("F9").Value = "#PBM:DWP:Notification Template:Investigated email template issue; found a rendering problem.  Template corrected.:R&D"
$ws.Range("F10").Value = "#PBM:DWP:Defect:DWP error caused by double quotes in request field; evaluate framework fix to handle special characters.:R&D"
$ws.Range("F11").Value = "#PBM:DWP:Branding:Confirmed separate logos for email template and DWP portal are possible.  Advised on customization steps.:Customization"
$ws.Range("F12").Value = "#PBM:DWP:Defect:UI refresh issue; corrected data display.:R&D"
$ws.Range("F13").Value = "#PBM:DWP:Request:Request to make DWP URLs public; feasibility assessment needed.:NA"
$ws.Range("F14").Value = "#PBM:DWP:Customization:Explore adding Excel support to Documents Form Element or suggest alternative solutions.:R&D"

Write-Host "Updated PBM_Tag values in F2:F14"
